$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'394"
$ws.Range("D2").Value = "'893828.79"
$ws.Range("C9").Value = "'121"
$ws.Range("D9").Value = "'267455.66"
$ws.Range("C10").Value = "'263"
$ws.Range("D10").Value = "'741381.77"
$ws.Range("C11").Value = "'112"
$ws.Range("D11").Value = "'284645.00"
$ws.Range("C12").Value = "'6"
$ws.Range("D12").Value = "'16250.00"
$ws.Range("C13").Value = "'169"
$ws.Range("D13").Value = "'405083.00"
$ws.Range("C15").Value = "'416"
$ws.Range("D15").Value = "'1334254.00"
$ws.Range("C16").Value = "'125"
$ws.Range("D16").Value = "'335901.38"
$ws.Range("C17").Value = "'11"
$ws.Range("D17").Value = "'26350.00"
$ws.Range("C38").Value = "'38"
$ws.Range("D38").Value = "'94956.00"
$ws.Range("C39").Value = "'36"
$ws.Range("D39").Value = "'154929.92"
$ws.Range("C40").Value = "'45"
$ws.Range("D40").Value = "'169599.15"
$ws.Range("C42").Value = "'210"
$ws.Range("D42").Value = "'538628.74"
$ws.Range("C44").Value = "'438"
$ws.Range("D44").Value = "'1397735.13"
$ws.Range("C45").Value = "'302"
$ws.Range("D45").Value = "'901358.86"
$ws.Range("C46").Value = "'5"
$ws.Range("D46").Value = "'10000.00"
$ws.Range("C47").Value = "'21"
$ws.Range("D47").Value = "'100220.65"
$ws.Range("C48").Value = "'1807"
$ws.Range("D48").Value = "'3750000.00"
$ws.Range("C50").Value = "'2789"
$ws.Range("D50").Value = "'7339449.03"
$ws.Range("C52").Value = "'2846"
$ws.Range("D52").Value = "'6862516.41"
$ws.Range("C66").Value = "'316"
$ws.Range("D66").Value = "'744618.84"
$ws.Range("C68").Value = "'784"
$ws.Range("D68").Value = "'2387232.47"
$ws.Range("C69").Value = "'453"
$ws.Range("D69").Value = "'1300292.91"
$ws.Range("C72").Value = "'355"
$ws.Range("D72").Value = "'764000.00"
$ws.Range("C74").Value = "'724"
$ws.Range("D74").Value = "'1952782.61"
$ws.Range("C75").Value = "'529"
$ws.Range("D75").Value = "'1292696.22"

Write-Output "done"
